# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 3 (pushing the existing rows 3..41 down
# to 4..42) and populate it with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 3..41 down to 4..42, growing the sheet to A1:R42.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new record.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44552
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112022
$ws.Range("G3").Value = "Arveja Verde"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 13500
$ws.Range("N3").Value = "$/saco 25 kilos"
$ws.Range("O3").Value = "Provincia de Diguillín"
$ws.Range("P3").Value = 540
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
